# Updates cryptocurrency Price (column D) and Volume(1h) (column E) values
# on the "cryptos" worksheet, per the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'23.918.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "

# Row 3
$ws.Range("D3").Value = "'1.655.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.76%  "

# Row 4
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.36%  "

# Row 5
$ws.Range("D5").Value = "'309.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.19%  "

# Row 6
$ws.Range("E6").Value = "  -0.26%  "

# Row 7
$ws.Range("D7").Value = "'0.3895"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.13%  "

# Row 8
$ws.Range("D8").Value = "'0.3841"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.88%  "

# Row 9
$ws.Range("D9").Value = "'51.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.86%  "

# Row 10
$ws.Range("D10").Value = "'1.354"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "

# Row 11
$ws.Range("D11").Value = "'1.003"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "

# Row 12
$ws.Range("D12").Value = "'0.08465"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.98%  "

# Row 13
$ws.Range("D13").Value = "'23.89"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.46%  "

# Row 14
$ws.Range("D14").Value = "'7.182"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.69%  "

# Row 15
$ws.Range("D15").Value = "'7.960"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.92%  "

# Row 16
$ws.Range("D16").Value = "'0.00001306"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.20%  "

# Row 17
$ws.Range("D17").Value = "'1.657.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.23%  "

# Row 18
$ws.Range("D18").Value = "'94.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.60%  "

# Row 19
$ws.Range("D19").Value = "'0.06990"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "

# Row 20
$ws.Range("D20").Value = "'19.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.46%  "

# Row 21
$ws.Range("D21").Value = "'6.930"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.47%  "

# Row 22
$ws.Range("E22").Value = "  -0.22%  "

# Row 23
$ws.Range("D23").Value = "'13.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.52%  "

# Row 24
$ws.Range("D24").Value = "'23.913.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.49%  "

# Row 25
$ws.Range("D25").Value = "'2.480"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.76%  "

# Row 26
$ws.Range("D26").Value = "'3.052"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.03%  "

# Row 27
$ws.Range("D27").Value = "'22.10"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.58%  "

# Row 28
$ws.Range("D28").Value = "'153.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "

# Row 29
$ws.Range("D29").Value = "'5.356"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.78%  "

# Row 30
$ws.Range("D30").Value = "'139.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.21%  "

# Row 31
$ws.Range("D31").Value = "'7.796"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.92%  "

# Row 32
$ws.Range("D32").Value = "'2.491"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.33%  "

# Row 33
$ws.Range("D33").Value = "'1.838.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.96%  "

# Row 34
$ws.Range("D34").Value = "'1.040"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.64%  "

# Row 35
$ws.Range("D35").Value = "'0.08075"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.96%  "

# Row 36
$ws.Range("D36").Value = "'0.02977"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.74%  "

# Row 37
$ws.Range("D37").Value = "'10.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.28%  "

# Row 38
$ws.Range("D38").Value = "'6.661"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.55%  "

# Row 39
$ws.Range("D39").Value = "'0.2690"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.96%  "

# Row 40
$ws.Range("D40").Value = "'0.09129"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.21%  "

# Row 41
$ws.Range("D41").Value = "'13.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.59%  "

# Row 42
$ws.Range("D42").Value = "'0.7531"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.87%  "

# Row 43
$ws.Range("D43").Value = "'1.413"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.19%  "

# Row 44
$ws.Range("D44").Value = "'16.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.82%  "

# Row 45
$ws.Range("D45").Value = "'0.6973"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.67%  "

# Row 46
$ws.Range("D46").Value = "'2.466"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.08%  "

# Row 47
$ws.Range("D47").Value = "'4.072"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

# Row 48
$ws.Range("E48").Value = "  -0.28%  "

# Row 49
$ws.Range("D49").Value = "'0.08263"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.66%  "

# Row 50
$ws.Range("D50").Value = "'134.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.24%  "

# Row 51
$ws.Range("D51").Value = "'1.233"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.19%  "
